$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, pushing existing rows 49-80 down to rows 50-81
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = 'Macroferia Regional de Talca'
$ws.Range("C49").Value = 'Maule'
$ws.Range("D49").Value = 44873
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = 300000000
$ws.Range("G49").Value = 'Espárragos'
$ws.Range("H49").Value = 'Sin especificar'
$ws.Range("I49").Value = 'Primera'
$ws.Range("J49").Value = 3000
$ws.Range("K49").Value = 1100
$ws.Range("L49").Value = 1100
$ws.Range("M49").Value = 1100
$ws.Range("N49").Value = '$/kilo'
$ws.Range("O49").Value = 'Provincia de Linares'
$ws.Range("P49").Value = 1100
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = 'Hortaliza'
